# Bill-of-materials update: swap a few parts for SMD equivalents, adjust
# their unit prices, and refresh the affected hyperlink targets.
#
# Commit message: "add sleep mode and assertions" (upstream message is
# generic / unrelated to the literal spreadsheet edit it accompanies —
# the actual content change is the BOM part/price/link refresh below).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Cell value edits (text + price updates)
# ---------------------------------------------------------------------

# Row 13: resistor renamed 300 -> 330 Ohm
$ws.Range("A13").Value = "Resistor 330 Ohm"

# Row 12: Resistor 10K -> cheaper SMD source
$ws.Range("D12").Value = 0.06
$ws.Range("F12").Value = "https://tpetrov.com/rezistor-10kom-025w-smd-1206-239372"

# Row 13: Resistor 330 Ohm -> SMD source
$ws.Range("F13").Value = "https://tpetrov.com/rezistor-330om-0125w-1-smd-243402"

# Row 14: Capacitor 2.2 uF -> cheaper SMD source
$ws.Range("D14").Value = 0.08
$ws.Range("F14").Value = "https://tpetrov.com/kondenzator-22uf-50v-smd-02-240100"

# Row 16: Capacitor 100 nF -> cheaper SMD source
$ws.Range("D16").Value = 0.03
$ws.Range("F16").Value = "https://tpetrov.com/kondenzator-100nf-25v-smd-0402-45886"

# ---------------------------------------------------------------------
# 2. Hyperlinks: F12 / F13 / F14 / F16 lose their hyperlink, F18 gains
#    a new one, and every other F-column hyperlink must survive with
#    its original "Hyperlink" cell style (fill/border) intact.
#
#    This runtime's Range.Hyperlinks.Delete()/.Item(n).Delete() cannot
#    remove a single link (only whole-sheet Hyperlinks.Delete() works),
#    and Hyperlinks.Add() always re-stamps the target cell with a fresh
#    "Hyperlink" style, dropping the sheet's custom fill/border. So:
#      a) snapshot each surviving link's cell (value+style) to a scratch
#         cell via Copy,
#      b) wipe all hyperlinks sheet-wide,
#      c) re-add only the links that should still exist,
#      d) copy each scratch snapshot back over its cell to restore the
#         original fill/border style (this does not disturb the
#         relationship Excel already associated with that cell address).
# ---------------------------------------------------------------------

# Cells whose hyperlink must be preserved, in the target rId order.
$keepOrder = @("F11", "F5", "F15", "F7", "F8", "F6", "F17", "F3")

# Cell that gains a brand-new hyperlink (same URL text it already shows).
$newLinkCell = "F18"
$newLinkUrl = "https://tpetrov.com/kondenzator-220mf-10v-105c-48326"

$scratchCols = @("H1", "H2", "H3", "H4", "H5", "H6", "H7", "H8", "H9")

# a) snapshot style+value of every surviving/linked cell before mutating
$idx = 0
foreach ($addr in $keepOrder) {
    $ws.Range($addr).Copy($ws.Range($scratchCols[$idx]))
    $idx = $idx + 1
}
$ws.Range($newLinkCell).Copy($ws.Range($scratchCols[$idx]))

# b) remove every hyperlink on the sheet (only working scope)
$ws.Hyperlinks.Delete()

# c) re-add the links that should exist, in order, so relationship ids
#    line up the same way the source workbook has them
foreach ($addr in $keepOrder) {
    $ws.Hyperlinks.Add($ws.Range($addr), $ws.Range($addr).Value())
}
$ws.Hyperlinks.Add($ws.Range($newLinkCell), $newLinkUrl)

# d) restore original fill/border style by copying the untouched
#    scratch snapshots back over their source cells
$idx = 0
foreach ($addr in $keepOrder) {
    $ws.Range($scratchCols[$idx]).Copy($ws.Range($addr))
    $idx = $idx + 1
}
$ws.Range($scratchCols[$idx]).Copy($ws.Range($newLinkCell))

# clean up scratch area
foreach ($c in $scratchCols) {
    $ws.Range($c).Clear()
}

# ---------------------------------------------------------------------
# 3. Selection moves to F18 (last thing the author clicked on)
# ---------------------------------------------------------------------
$ws.Range("F18").Select()
